# Applies the "Perfil usuario, terminos y condiciones" edit:
#  1. Merge the "Configure Shopify admin..." + "client" runs (drop stray proofErr marks)
#  2. Strike through "Build Shopping Cart component" (keep " and page" un-struck)
#  3. Strike through "Implement product search functionality"
#  4. Strike through the "(returns, privacy, terms) in Shopify" part only
#  5. Strike through 'Build User Profile page ("Perfil Usuario")'
#  6. Append a new bold "A cambiar antes de la entrega:" note paragraph at the end

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# --- Change 1: "Configure Shopify admin and create documentation for client" ---
$i1 = Find-ParagraphIndex $d "Configure Shopify admin and create documentation for"
$d.Paragraphs($i1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Configure Shopify admin and create documentation for client</w:t></w:r></w:p>')

# --- Change 2: "Build Shopping Cart component and page" ---
$i2 = Find-ParagraphIndex $d "Build Shopping Cart component and page"
$d.Paragraphs($i2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Build Shopping Cart component</w:t></w:r><w:r><w:t xml:space="preserve"> and page</w:t></w:r></w:p>')

# --- Change 3: "Implement product search functionality" ---
$i3 = Find-ParagraphIndex $d "Implement product search functionality"
$d.Paragraphs($i3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:lastRenderedPageBreak/><w:t>Implement product search functionality</w:t></w:r></w:p>')

# --- Change 4: "Create policies pages (returns, privacy, terms) in Shopify" ---
$i4 = Find-ParagraphIndex $d "Create policies pages (returns, privacy, terms) in Shopify"
$d.Paragraphs($i4).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Create policies pages </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>(returns, privacy, terms) in Shopify</w:t></w:r></w:p>')

# --- Change 5: Build User Profile page ("Perfil Usuario") ---
$i5 = Find-ParagraphIndex $d 'Build User Profile page ("Perfil Usuario")'
$d.Paragraphs($i5).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Build User Profile page ("Perfil Usuario")</w:t></w:r></w:p>')

# --- Change 6: append the new note paragraph after the last paragraph in the doc ---
$iLast = $d.Paragraphs.Count
$d.Paragraphs($iLast).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:strike/></w:rPr><w:t>The 3D dog feature</w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> might need scoping - is it just a viewer or interactive configurator?</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>cambiar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> antes de la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>entrega</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:r><w:t>- Correo store owner Shopify</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Agregar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dominio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> nuevo a callbacks Headless&gt;Customer API</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Vercel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Envs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Next </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>urls</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(2)</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sho</w:t></w:r><w:r><w:t>pify</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> customer callback </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, store owner email, </w:t></w:r></w:p>')

Write-Output "All changes applied"
